$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetime for the 653de435... file (row 2)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-18 09:54:23"
$wsZh.Range("G2").Value = "2016-02-18 09:55:08"

# de-de sheet: update handoff/handback datetime for the 653de435... file (row 2)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-18 09:54:33"
$wsDe.Range("G2").Value = "2016-02-18 09:55:31"
